# Fix microsoft node link: update onnx node ids referenced inside the
# error-message strings on the "Training Results" sheet (sheet1 / rId1).
#
# Each cell below contains a long inline-string error message that embeds a
# numeric onnx node id (e.g. onnx::Pow::27433). The node ids changed after a
# re-run against a newer onnxruntime build, so only the numeric id inside the
# text needs to be updated; everything else in the string stays identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Results")

$ws.Range("C16").Value = "C:\Users\COCO\onnxruntime_training_cuda_python\orttraining\orttraining\python\orttraining_pybind_state.cc:621 onnxruntime::python::addObjectMethodsForTraining::<lambda_6dd399ad6691adab5d0e0423ed8ce22d>::operator () [ONNXRuntimeError] : 1 : FAIL : Type Error: Type parameter (T) of Optype (Sub) bound to different types (tensor(float) and tensor(double) in node (onnx::Pow::32375_Grad/Sub_1).`n"

$ws.Range("C39").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32494): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C51").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32584): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C52").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32586): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C64").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32639): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C65").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32641): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C69").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32659): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C70").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32661): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C93").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::32757): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C102").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::32794): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C223").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::33244): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C239").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::33249): X typestr: T, has unsupported type: tensor(uint8)"
